$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "video crashed" -> "video crashed :(" (Arrays - Problem Solving row, column F)
$ws.Range("F9").Value2 = "video crashed :("

# 2. Fix the copy/pasted lesson numbers for the OOP section (rows 18-20 were all
#    left at 14 after the row above them; bump them to the correct sequence).
$ws.Range("B18").Value2 = 15
$ws.Range("B19").Value2 = 16
$ws.Range("B20").Value2 = 17

# 3. Add the recording dates for "OOP #2" and "OOP #3", matching the date format
#    already used throughout column E.
$ws.Range("E18").NumberFormat = $ws.Range("E17").NumberFormat
$ws.Range("E18").Value2 = 44151
$ws.Range("E19").NumberFormat = $ws.Range("E17").NumberFormat
$ws.Range("E19").Value2 = 44156

# 4. Add placeholder link text for "OOP #2" (slides attached separately).
$ws.Range("F18").Value2 = "Part #1: Part #2:"

# 5. Move selection / scroll position, mirroring where the editor left off.
$ws.Range("E20").Select()
